$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49, pushing the existing rows 49:62 down to 50:63
$ws.Rows("49:49").Insert()

# Populate the newly inserted row 49 with the new record
$ws.Cells.Item(49, 1).Value = 1
$ws.Cells.Item(49, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(49, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(49, 4).Value = 44468
$ws.Cells.Item(49, 5).Value = 15
$ws.Cells.Item(49, 6).Value = "Fruta"
$ws.Cells.Item(49, 7).Value = 100106
$ws.Cells.Item(49, 8).Value = "Oleaginosos"
$ws.Cells.Item(49, 9).Value = 100106002
$ws.Cells.Item(49, 10).Value = "Palta"
$ws.Cells.Item(49, 11).Value = "Edranol"
$ws.Cells.Item(49, 12).Value = "Tercera"
$ws.Cells.Item(49, 13).Value = 250
$ws.Cells.Item(49, 14).Value = 54000
$ws.Cells.Item(49, 15).Value = 55000
$ws.Cells.Item(49, 16).Value = 54500
$ws.Cells.Item(49, 17).Value = "$/caja 25 kilos"
$ws.Cells.Item(49, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(49, 19).Value = 2180
$ws.Cells.Item(49, 20).Value = 25
